# Require angle brackets for @base / @prefix IRI values on the
# "North Wind Info" sheet.
#
# Column layout (B: @prefix / C: prefix-label / D: IRI value), rows 1-6
# are the @prefix declarations, row 7 is the rdfs:subclassOf triple.
#
# Statement order matters: new literal strings are appended to the shared
# strings table in the order they are first written, and the target
# workbook expects:
#   ... :data, :schema,
#   <http://seman.tc/data/northwind>,
#   <http://seman.tc/models/northwind#>,
#   <http://schema.org/>,
#   <http://xmlns.com/foaf/0.1/>,
#   <http://purl.org/dc/terms/>,
#   <http://seman.tc/data/northwind/ApplicationModule>,
#   <http://seman.tc/data/northwind/VCAMPApplicationModule>
# so the D-column (IRI) cells are written before the B7 (VCAMP module IRI)
# cell, and D7 is written ahead of B7 to land ApplicationModule before
# VCAMPApplicationModule.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("North Wind Info")

# Row 1 : @prefix :data <http://seman.tc/data/northwind>
$ws.Range("D1").Value = "<http://seman.tc/data/northwind>"

# Row 2 : @prefix :schema <http://seman.tc/models/northwind#>
$ws.Range("D2").Value = "<http://seman.tc/models/northwind#>"

# Row 3 : @prefix schema: <http://schema.org/>
$ws.Range("C3").Value = "schema:"
$ws.Range("D3").Value = "<http://schema.org/>"

# Row 4 : @prefix foaf: <http://xmlns.com/foaf/0.1/>
$ws.Range("C4").Value = "foaf:"
$ws.Range("D4").Value = "<http://xmlns.com/foaf/0.1/>"

# Row 5 : @prefix northwind: <http://seman.tc/models/northwind#>
$ws.Range("C5").Value = "northwind:"
$ws.Range("D5").Value = "<http://seman.tc/models/northwind#>"

# Row 6 : @prefix dct: <http://purl.org/dc/terms/>
$ws.Range("C6").Value = "dct:"
$ws.Range("D6").Value = "<http://purl.org/dc/terms/>"

# Row 7 : <VCAMPApplicationModule> rdfs:subclassOf <ApplicationModule>
$ws.Range("D7").Value = "<http://seman.tc/data/northwind/ApplicationModule>"
$ws.Range("B7").Value = "<http://seman.tc/data/northwind/VCAMPApplicationModule>"

# Re-fit columns B:D to the new (longer) content, mirroring Excel's
# AutoFit-on-edit behaviour, then leave the selection on B8 (below the
# edited table) as recorded in the saved workbook.
$ws.Columns.Item(2).ColumnWidth = 55.0
$ws.Columns.Item(3).ColumnWidth = 13.666666666666666
$ws.Columns.Item(4).ColumnWidth = 50.333333333333336

$ws.Range("B8").Select()
